# Add the "accent3" theme color (9BBB59) to the paragraph-mark run
# properties (pPr/rPr) and the run properties (r/rPr) of the single
# paragraph that just holds the YouTube link
# "https://www.youtube.com/watch?v=f4RR5YO391w" -- matching the color
# already applied to the sibling link paragraphs right above it.

$d = $word.ActiveDocument

$target = "f4RR5YO391w"

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*$target*") {

        # Rebuild this paragraph's OOXML with <w:color w:val="9BBB59"
        # w:themeColor="accent3"/> inserted into both the paragraph
        # mark's run properties and the run's run properties, right
        # after <w:rFonts .../> -- everything else (rsids, fonts,
        # size, language) stays exactly as it was.
        $fragment = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="006221F9" w:rsidRDefault="00F3030B"><w:pPr><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="9BBB59" w:themeColor="accent3"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="fr-FR"/></w:rPr></w:pPr><w:r w:rsidRPr="00F3030B"><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:color w:val="9BBB59" w:themeColor="accent3"/><w:sz w:val="26"/><w:szCs w:val="26"/><w:lang w:val="fr-FR"/></w:rPr><w:t>https://www.youtube.com/watch?v=f4RR5YO391w</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

        $p.Range.InsertXML($fragment)
        break
    }
}
